# ----------------------------------------------------------------------
# Edit: Thu, Jun 25, 2020 10:04:51 PM
#
# 1) The table on the "Component three" slide switches from the custom
#    "Table_0" style to the built-in PowerPoint table style
#    {225D3552-7D1D-411F-A096-8E113EDA1A11}.
# 2) The presentation's theme (color scheme) is changed from the
#    custom "Integral" palette to the stock "Office" palette (i.e. the
#    Design tab's default "Office Theme" was (re)applied).
# ----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table -------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{225D3552-7D1D-411F-A096-8E113EDA1A11}")
        }
    }
}

# --- 2. Swap the theme colours from "Integral" to the default "Office" ---
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Item(1).RGB  = 0          # Dark 1       000000
$themeColors.Item(2).RGB  = 16777215   # Light 1      FFFFFF
$themeColors.Item(3).RGB  = 6968388    # Dark 2       44546A
$themeColors.Item(4).RGB  = 15132391   # Light 2      E7E6E6
$themeColors.Item(5).RGB  = 13998939   # Accent 1     5B9BD5
$themeColors.Item(6).RGB  = 3243501    # Accent 2     ED7D31
$themeColors.Item(7).RGB  = 10855845   # Accent 3     A5A5A5
$themeColors.Item(8).RGB  = 49407      # Accent 4     FFC000
$themeColors.Item(9).RGB  = 12874308   # Accent 5     4472C4
$themeColors.Item(10).RGB = 4697456    # Accent 6     70AD47
$themeColors.Item(11).RGB = 12673797   # Hyperlink    0563C1
$themeColors.Item(12).RGB = 7491477    # Followed Hyperlink 954F72
